$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the old row 27 (the last, thick-bottom-bordered
# data row). This pushes old row 27 down to row 29, and the totals /
# notes / footer rows below it shift down by 2 as well.
$ws.Rows("27:28").Insert(-4121, 0)

# The generic Insert() call above does not reliably carry over the
# "regular" (non-thick-bottom) row formatting from row 26, so copy it
# explicitly onto the two freshly inserted rows.
$ws.Range("B26:K26").Copy()
$ws.Range("B27:K28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows("27:28").RowHeight = $ws.Rows("26:26").RowHeight

# Fill in the row-number column (A) for the newly inserted data rows,
# and renumber the old row (now pushed down to row 29) from 23 to 25.
$ws.Range("A27").Value = 23
$ws.Range("A28").Value = 24
$ws.Range("A29").Value = 25

# Update the totals row formulas (now row 30) to include the two new
# data rows (the SUM ranges grow from D5:D27 to D5:D29, etc.)
$ws.Range("D30").Formula = "=SUM(D5:D29)"
$ws.Range("E30").Formula = "=SUM(E5:E29)"
$ws.Range("F30").Formula = "=SUM(F5:F29)"
$ws.Range("G30").Formula = "=SUM(G5:G29)"
$ws.Range("H30").Formula = "=SUM(H5:H29)"
$ws.Range("I30").Formula = "=F30/H30"
$ws.Range("J30").Formula = "=F30/G30"
$ws.Range("K30").Formula = "=G30/H30"

# Restore the print area (grew by 2 rows) and the sheet-view selection.
$ws.PageSetup.PrintArea = "`$B`$1:`$K`$33"
$ws.Range("B2:B4").Select() | Out-Null
